# Generate Report for Handback
# Update the localization-status report to reflect that the d6875cf9-... file
# has now been handed back (in sync with en-US), instead of merely being
# "Ready for handoff" / showing a stale-handback error.

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

# --- Overview sheet: the zh-cn / de-de status columns for the
#     d6875cf9-4019-4dc6-b293-0d94272724f5.md row (row 3) ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus

# --- zh-cn sheet: Status / Latest Handback DateTime / Error Detail for
#     the d6875cf9-... row (row 3) ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = $newStatus
$wsZhCn.Range("K3").Value = "2016-08-19 18:56:29"
$wsZhCn.Range("P3").Value = ""

# --- de-de sheet: Status / Latest Handback DateTime / Error Detail for
#     the d6875cf9-... row (row 3) ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = $newStatus
$wsDeDe.Range("K3").Value = "2016-08-19 18:56:36"
$wsDeDe.Range("P3").Value = ""
